$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2851076882107093
$ws.Range("C2").Value = 2.671180954222684
$ws.Range("B3").Value = 0.2991427698894775
$ws.Range("C3").Value = 3.086674271150811
